# Add a new "2022" column (column O) to the report, mirroring the existing
# year columns (D..N = 2011..2021).
#
# Styling plan (matches the target cellXfs layout):
#   - O2  : blank cell, same formatting as N2  (existing style, no border side cell)
#   - O3  : header value 2022,                  same formatting as N3 (existing style)
#   - O4  : bold data row,                      new style (numFmt 0.0, bold 9pt Times New Roman, vertical-center, no border)
#   - O5..O14 : regular data rows,              new style (numFmt 0.0, regular 9pt Times New Roman, vertical-center, no border)
#   - O15 : footer value 100,                   same formatting as N15 (existing style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- O2: empty cell, formatted like N2 -------------------------------------
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# --- O3: year header, formatted like N3 -------------------------------------
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = 2022

# --- O4: bold data row (first data row gets its own bold style) -------------
$ws.Range("O4").Value = 96.345513960706299
$ws.Range("O4").NumberFormat = "0.0"
$ws.Range("O4").Font.Name = "Times New Roman"
$ws.Range("O4").Font.Size = 9
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").VerticalAlignment = -4108

# --- O5: first regular data row establishes the shared non-bold style -------
$ws.Range("O5").Value = 99.646905185978142
$ws.Range("O5").NumberFormat = "0.0"
$ws.Range("O5").Font.Name = "Times New Roman"
$ws.Range("O5").Font.Size = 9
$ws.Range("O5").VerticalAlignment = -4108

# --- O6:O14: same regular style as O5, with their own values ----------------
$ws.Range("O5").Copy()
$ws.Range("O6:O14").PasteSpecial(-4122)

$ws.Range("O6").Value = 94.429993069436605
$ws.Range("O7").Value = 88.286387066773813
$ws.Range("O8").Value = 93.885244420521602
$ws.Range("O9").Value = 98.153167726175582
$ws.Range("O10").Value = 97.951019527503291
$ws.Range("O11").Value = 93.640014938442292
$ws.Range("O12").Value = 99.643271453928278
$ws.Range("O13").Value = 100
$ws.Range("O14").Value = 100

# --- O15: footer total, formatted like N15 -----------------------------------
$ws.Range("N15").Copy()
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("O15").Value = 100

# --- Update the active selection to mirror the shifted "next empty column" --
$ws.Range("P2").Select()
